# Auto commit: 2025-04-15 21:43:49
# Rename "MinVarPortfolio" -> "MVP_Weights" and add a new "MVP_Stats" sheet
# with portfolio mean return / standard deviation statistics.

$wb = $excel.ActiveWorkbook

# 1) Rename the existing MinVarPortfolio sheet to MVP_Weights
$mvp = $wb.Worksheets.Item("MinVarPortfolio")
$mvp.Name = "MVP_Weights"

# 2) Add a brand-new worksheet after the last existing sheet, named MVP_Stats
$statsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$statsSheet.Name = "MVP_Stats"

# 3) Populate the new sheet with the portfolio statistics
$statsSheet.Range("B1").Value = "Portfolio Mean Return"
$statsSheet.Range("C1").Value = "Portfolio Standard Deviation"
$statsSheet.Range("A2").Value = "Portfolio"
$statsSheet.Range("B2").Value = 0.00108357142866864
$statsSheet.Range("C2").Value = 0.01091700535987142

# 4) Match the header styling used elsewhere in the workbook (bold, centered,
#    bordered) by copying the format from an existing styled header cell.
$headerStyleSource = $wb.Worksheets.Item("Statistics").Range("A1")

$headerStyleSource.Copy()
$statsSheet.Range("B1:C1").PasteSpecial(-4122)

$headerStyleSource.Copy()
$statsSheet.Range("A2").PasteSpecial(-4122)
